# IPDO.xlsx - "imprimir_resultados" routine now wraps its logging in a
# try/catch; when it runs again it appends a fresh snapshot row to the
# BalancoResumido log sheet instead of overwriting the previous one, so the
# sheet ends up with two more rows than before (the duplicated log entry),
# and the previously-last row's trailing marker cell moves to the new
# last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BalancoResumido")

# Append the latest log snapshot twice (rows 86 and 87), duplicating the
# values that were already on the last existing row (85).
$ws.Range("A85:Q85").Copy($ws.Range("A86:Q86"))
$ws.Range("A85:R85").Copy($ws.Range("A87:R87"))

# The trailing marker cell that used to sit at the end of the old last row
# now belongs to the new last row only.
$ws.Range("R85").ClearContents()
